$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# BF2:BF31 hold the game-date label for each team row. The source data was
# off by one day because of how the NBA stats site reported the date, so
# every "5-10-2011-12" becomes the corrected ISO-style "2012-05-10".
$dateRange = $ws.Range("BF2:BF31")

# Force text storage first so Excel's automatic date-literal detection
# doesn't silently convert the assigned string into a date serial number.
$dateRange.NumberFormat = "@"

for ($row = 2; $row -le 31; $row++) {
    $ws.Cells.Item($row, 58).Value = "2012-05-10"
}

# Restore the default (unstyled) cell style now that the text is safely
# stored as a string, matching the original unformatted BF2:BF31 cells.
$dateRange.Style = "Normal"
